$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.3674653333333333
$ws.Range("H2").Value = 1.102396
$ws.Range("I2").Value = 0.5572612813719676
$ws.Range("J2").Value = 0.5572612813719677
$ws.Range("M2").Value = 45.1830845
$ws.Range("N2").Value = 90.366169
$ws.Range("O2").Value = 0.2982772948921854
$ws.Range("P2").Value = 0.2359735829156887
$ws.Range("Q2").Value = 16.60321720682067
$ws.Range("R2").Value = 99.619303240924
$ws.Range("S2").Value = 0.1662183875557835
$ws.Range("T2").Value = 0.131498941185531
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.3674653333333333
$ws.Range("H3").Value = 1.102396
$ws.Range("I3").Value = 0.5572612813719676
$ws.Range("J3").Value = 0.5572612813719677
$ws.Range("N3").Value = 73.46982600000001
$ws.Range("O3").Value = 0.1616709822417395
$ws.Range("P3").Value = 0.1918520865636367
$ws.Range("Q3").Value = 8.999204700344
$ws.Range("R3").Value = 80.99284230309601
$ws.Range("S3").Value = 0.09009297872469639
$ws.Range("T3").Value = 0.1069117395923379
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.3674653333333333
$ws.Range("H4").Value = 1.102396
$ws.Range("I4").Value = 0.5572612813719676
$ws.Range("J4").Value = 0.5572612813719677
$ws.Range("M4").Value = 20.755341
$ws.Range("N4").Value = 62.26602299999999
$ws.Range("O4").Value = 0.1370169176485697
$ws.Range("P4").Value = 0.1625955454769879
$ws.Range("Q4").Value = 7.626868299011999
$ws.Range("R4").Value = 68.64181469110798
$ws.Range("S4").Value = 0.07635422309847931
$ws.Range("T4").Value = 0.09060820201788035
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.3674653333333333
$ws.Range("H5").Value = 1.102396
$ws.Range("I5").Value = 0.5572612813719676
$ws.Range("J5").Value = 0.5572612813719677
$ws.Range("M5").Value = 26.3069545
$ws.Range("N5").Value = 52.613909
$ws.Range("O5").Value = 0.173666037012409
$ws.Range("P5").Value = 0.1373909368441856
$ws.Range("Q5").Value = 9.666893804327332
$ws.Range("R5").Value = 58.00136282596399
$ws.Range("S5").Value = 0.0967773583163266
$ws.Range("T5").Value = 0.07656264951468594
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.3674653333333333
$ws.Range("H6").Value = 1.102396
$ws.Range("I6").Value = 0.5572612813719676
$ws.Range("J6").Value = 0.5572612813719677
$ws.Range("M6").Value = 17.34473466666667
$ws.Range("N6").Value = 52.034204
$ws.Range("O6").Value = 0.1145017121838161
$ws.Range("P6").Value = 0.1358771505744131
$ws.Range("Q6").Value = 6.373588705864889
$ws.Range("R6").Value = 57.362298352784
$ws.Range("S6").Value = 0.06380737085083762
$ws.Range("T6").Value = 0.07571907503826925
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.3674653333333333
$ws.Range("H7").Value = 1.102396
$ws.Range("I7").Value = 0.5572612813719676
$ws.Range("J7").Value = 0.5572612813719677
$ws.Range("M7").Value = 17.400077
$ws.Range("N7").Value = 52.200231
$ws.Range("O7").Value = 0.1148670560212801
$ws.Range("P7").Value = 0.136310697625088
$ws.Range("Q7").Value = 6.393925094830666
$ws.Range("R7").Value = 57.54532585347599
$ws.Range("S7").Value = 0.06401096282584413
$ws.Range("T7").Value = 0.07596067402326338
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.2919476666666667
$ws.Range("H8").Value = 0.875843
$ws.Range("I8").Value = 0.4427387186280323
$ws.Range("J8").Value = 0.4427387186280324
$ws.Range("M8").Value = 45.1830845
$ws.Range("N8").Value = 90.366169
$ws.Range("O8").Value = 0.2982772948921854
$ws.Range("P8").Value = 0.2359735829156887
$ws.Range("Q8").Value = 13.19109609257783
$ws.Range("R8").Value = 79.14657655546701
$ws.Range("S8").Value = 0.1320589073364019
$ws.Range("T8").Value = 0.1044746417301578
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.2919476666666667
$ws.Range("H9").Value = 0.875843
$ws.Range("I9").Value = 0.4427387186280323
$ws.Range("J9").Value = 0.4427387186280324
$ws.Range("N9").Value = 73.46982600000001
$ws.Range("O9").Value = 0.1616709822417395
$ws.Range("P9").Value = 0.1918520865636367
$ws.Range("Q9").Value = 7.149781423702001
$ws.Range("R9").Value = 64.34803281331801
$ws.Range("S9").Value = 0.07157800351704312
$ws.Range("T9").Value = 0.08494034697129887
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.2919476666666667
$ws.Range("H10").Value = 0.875843
$ws.Range("I10").Value = 0.4427387186280323
$ws.Range("J10").Value = 0.4427387186280324
$ws.Range("M10").Value = 20.755341
$ws.Range("N10").Value = 62.26602299999999
$ws.Range("O10").Value = 0.1370169176485697
$ws.Range("P10").Value = 0.1625955454769879
$ws.Range("Q10").Value = 6.059473375820999
$ws.Range("R10").Value = 54.53526038238899
$ws.Range("S10").Value = 0.06066269455009037
$ws.Range("T10").Value = 0.07198734345910761
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.2919476666666667
$ws.Range("H11").Value = 0.875843
$ws.Range("I11").Value = 0.4427387186280323
$ws.Range("J11").Value = 0.4427387186280324
$ws.Range("M11").Value = 26.3069545
$ws.Range("N11").Value = 52.613909
$ws.Range("O11").Value = 0.173666037012409
$ws.Range("P11").Value = 0.1373909368441856
$ws.Range("Q11").Value = 7.680253983381166
$ws.Range("R11").Value = 46.081523900287
$ws.Range("S11").Value = 0.07688867869608239
$ws.Range("T11").Value = 0.06082828732949964
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.2919476666666667
$ws.Range("H12").Value = 0.875843
$ws.Range("I12").Value = 0.4427387186280323
$ws.Range("J12").Value = 0.4427387186280324
$ws.Range("M12").Value = 17.34473466666667
$ws.Range("N12").Value = 52.034204
$ws.Range("O12").Value = 0.1145017121838161
$ws.Range("P12").Value = 0.1358771505744131
$ws.Range("Q12").Value = 5.063754814885778
$ws.Range("R12").Value = 45.57379333397201
$ws.Range("S12").Value = 0.05069434133297851
$ws.Range("T12").Value = 0.06015807553614387
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.2919476666666667
$ws.Range("H13").Value = 0.875843
$ws.Range("I13").Value = 0.4427387186280323
$ws.Range("J13").Value = 0.4427387186280324
$ws.Range("M13").Value = 17.400077
$ws.Range("N13").Value = 52.200231
$ws.Range("O13").Value = 0.1148670560212801
$ws.Range("P13").Value = 0.136310697625088
$ws.Range("Q13").Value = 5.079911879970333
$ws.Range("R13").Value = 45.719206919733
$ws.Range("S13").Value = 0.05085609319543595
$ws.Range("T13").Value = 0.06035002360182465
